# Updates cryptos list (prices/volumes) and re-orders three coin pairs,
# matching the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are text (e.g. "42.127.92", "0.630") - force text format
# so Excel does not reinterpret them as numbers/dates and strip formatting.

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '42.127.92'
$ws.Range('E2').Value = '  -2.31%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '2.233.44'
$ws.Range('E3').Value = '  -3.04%  '
$ws.Range('E4').Value = '  +0.19%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '246.05'
$ws.Range('E5').Value = '  -2.76%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '0.630'
$ws.Range('E6').Value = '  -2.17%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '76.12'
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  +0.03%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.617'
$ws.Range('E9').Value = '  -5.28%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '41.66'
$ws.Range('E10').Value = '  +5.59%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.0942'
$ws.Range('E11').Value = '  -4.75%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '7.03'
$ws.Range('E12').Value = '  -8.31%  '
$ws.Range('E13').Value = '  -3.69%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '2.568.40'
$ws.Range('E14').Value = '  -3.07%  '
$ws.Range('E15').Value = '  -5.68%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '0.851'
$ws.Range('E16').Value = '  -3.35%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '2.239.30'
$ws.Range('E17').Value = '  -2.85%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '41.939.49'
$ws.Range('E18').Value = '  -2.59%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0978'
$ws.Range('E19').Value = '  -3.75%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '71.54'
$ws.Range('E20').Value = '  -2.02%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '6.05'
$ws.Range('E21').Value = '  -4.23%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '2.27'
$ws.Range('E22').Value = '  +0.37%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '230.48'
$ws.Range('E23').Value = '  -3.42%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$ws.Range('E24').Value = '  -0.14%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '11.21'
$ws.Range('E25').Value = '  -4.07%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '3.68'
$ws.Range('E26').Value = '  -5.78%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '2.28'
$ws.Range('E27').Value = '  -5.84%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '7.34'
$ws.Range('E28').Value = '  +14.77%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '2.15'
$ws.Range('E29').Value = '  +0.33%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '169.03'
$ws.Range('E30').Value = '  +0.55%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '20.46'
$ws.Range('E31').Value = '  -3.79%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '0.0822'
$ws.Range('E32').Value = '  -3.08%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '32.16'
$ws.Range('E33').Value = '  +4.32%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '0.119'
$ws.Range('E34').Value = '  -6.36%  '
$ws.Range('E35').Value = '  -2.39%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '4.43'
$ws.Range('E36').Value = '  -3.89%  '
$ws.Range('B37').Value = 'Filecoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '4.91'
$ws.Range('E37').Value = '  +0.72%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.0299'
$ws.Range('E38').Value = '  -4.60%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '13.95'
$ws.Range('E39').Value = '  +0.86%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '5.84'
$ws.Range('E40').Value = '  -0.56%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '2.16'
$ws.Range('E41').Value = '  -8.57%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '111.66'
$ws.Range('E42').Value = '  +6.50%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '0.201'
$ws.Range('E43').Value = '  -8.63%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '60.41'
$ws.Range('E44').Value = '  -3.64%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '8.63'
$ws.Range('E45').Value = '  -6.65%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '0.0986'
$ws.Range('E46').Value = '  -5.04%  '
$ws.Range('B47').Value = 'BinanceUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '0.997'
$ws.Range('E47').Value = '  -0.57%  '
$ws.Range('E48').Value = '  -5.33%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '1.16'
$ws.Range('E49').Value = '  -2.27%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '4.20'
$ws.Range('E50').Value = '  -13.87%  '
$ws.Range('E51').Value = '  +12.31%  '
